# Timesheet changes by Ruchika
# Fill in the Feb-18 (Saturday, column X) "OFF" marker and the actual
# hours worked for Feb-16/17/19 (columns V/W/Y) on rows 28-31 of the
# "February 2013" sheet - mirroring the existing OFF/weekend pattern
# already present in column Q (Feb-11, the previous Saturday).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("February 2013")

# Column X (Saturday 18-Feb) gets the same "OFF" styling/text already used
# for column Q (Saturday 11-Feb) on these rows - copy the format over first
# so the OFF shading comes across, then stamp in the text.
$ws.Range("Q28").Copy()
$ws.Range("X28").PasteSpecial(-4122)
$ws.Range("Q29").Copy()
$ws.Range("X29").PasteSpecial(-4122)
$ws.Range("Q30").Copy()
$ws.Range("X30").PasteSpecial(-4122)
$ws.Range("Q31").Copy()
$ws.Range("X31").PasteSpecial(-4122)

$ws.Range("X28").Value = "OFF"
$ws.Range("X29").Value = "OFF"
$ws.Range("X30").Value = "OFF"
$ws.Range("X31").Value = "OFF"

# Hours logged for the surrounding days.
$ws.Range("V28").Value = 0
$ws.Range("W28").Value = 0
$ws.Range("Y28").Value = 0

$ws.Range("V29").Value = 0
$ws.Range("W29").Value = 0
$ws.Range("Y29").Value = 0

$ws.Range("V30").Value = 2
$ws.Range("W30").Value = 1
$ws.Range("Y30").Value = 2

$ws.Range("V31").Value = 0
$ws.Range("W31").Value = 0.5
$ws.Range("Y31").Value = 2

# Restore the active selection/scroll position recorded for the sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 23
$win.ScrollColumn = 18
$ws.Range("Z34").Select() | Out-Null
